# Natmi following Dr Hou advice
# Rebuilds the LR-pair results table for Artn-Gfra1 so that every
# combination of Sending cluster {FAPs, sCs} x Target cluster {ECs, FAPs, sCs}
# (Ligand = Artn, Receptor = Gfra1) is present, rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Artn/Gfra1 -> ECs
$ws.Cells.Item(2, 1).Value  = "FAPs"
$ws.Cells.Item(2, 2).Value  = "Artn"
$ws.Cells.Item(2, 3).Value  = "Gfra1"
$ws.Cells.Item(2, 4).Value  = "ECs"
$ws.Cells.Item(2, 5).Value  = 3
$ws.Cells.Item(2, 6).Value  = 1
$ws.Cells.Item(2, 7).Value  = 1.415594
$ws.Cells.Item(2, 8).Value  = 4.246782
$ws.Cells.Item(2, 9).Value  = 0.9753967820130081
$ws.Cells.Item(2, 10).Value = 0.9753967820130081
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.509196
$ws.Cells.Item(2, 14).Value = 1.527588
$ws.Cells.Item(2, 15).Value = 0.02558190413389134
$ws.Cells.Item(2, 16).Value = 0.02558190413389134
$ws.Cells.Item(2, 17).Value = 0.7208148024239999
$ws.Cells.Item(2, 18).Value = 6.487333221815999
$ws.Cells.Item(2, 19).Value = 0.02495250696996288
$ws.Cells.Item(2, 20).Value = 0.02495250696996288

# Row 3: FAPs -> Artn/Gfra1 -> FAPs
$ws.Cells.Item(3, 1).Value  = "FAPs"
$ws.Cells.Item(3, 2).Value  = "Artn"
$ws.Cells.Item(3, 3).Value  = "Gfra1"
$ws.Cells.Item(3, 4).Value  = "FAPs"
$ws.Cells.Item(3, 5).Value  = 3
$ws.Cells.Item(3, 6).Value  = 1
$ws.Cells.Item(3, 7).Value  = 1.415594
$ws.Cells.Item(3, 8).Value  = 4.246782
$ws.Cells.Item(3, 9).Value  = 0.9753967820130081
$ws.Cells.Item(3, 10).Value = 0.9753967820130081
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 17.676258
$ws.Cells.Item(3, 14).Value = 53.028774
$ws.Cells.Item(3, 15).Value = 0.8880516296316739
$ws.Cells.Item(3, 16).Value = 0.8880516296316739
$ws.Cells.Item(3, 17).Value = 25.022404767252
$ws.Cells.Item(3, 18).Value = 225.201642905268
$ws.Cells.Item(3, 19).Value = 0.8662027018041424
$ws.Cells.Item(3, 20).Value = 0.8662027018041424

# Row 4: FAPs -> Artn/Gfra1 -> sCs
$ws.Cells.Item(4, 1).Value  = "FAPs"
$ws.Cells.Item(4, 2).Value  = "Artn"
$ws.Cells.Item(4, 3).Value  = "Gfra1"
$ws.Cells.Item(4, 4).Value  = "sCs"
$ws.Cells.Item(4, 5).Value  = 3
$ws.Cells.Item(4, 6).Value  = 1
$ws.Cells.Item(4, 7).Value  = 1.415594
$ws.Cells.Item(4, 8).Value  = 4.246782
$ws.Cells.Item(4, 9).Value  = 0.9753967820130081
$ws.Cells.Item(4, 10).Value = 0.9753967820130081
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.719084666666667
$ws.Cells.Item(4, 14).Value = 5.157254
$ws.Cells.Item(4, 15).Value = 0.0863664662344347
$ws.Cells.Item(4, 16).Value = 0.0863664662344347
$ws.Cells.Item(4, 17).Value = 2.433525939625333
$ws.Cells.Item(4, 18).Value = 21.901733456628
$ws.Cells.Item(4, 19).Value = 0.08424157323890272
$ws.Cells.Item(4, 20).Value = 0.08424157323890272

# Row 5: sCs -> Artn/Gfra1 -> ECs
$ws.Cells.Item(5, 1).Value  = "sCs"
$ws.Cells.Item(5, 2).Value  = "Artn"
$ws.Cells.Item(5, 3).Value  = "Gfra1"
$ws.Cells.Item(5, 4).Value  = "ECs"
$ws.Cells.Item(5, 5).Value  = 1
$ws.Cells.Item(5, 6).Value  = 0.3333333333333333
$ws.Cells.Item(5, 7).Value  = 0.03570666666666667
$ws.Cells.Item(5, 8).Value  = 0.10712
$ws.Cells.Item(5, 9).Value  = 0.0246032179869919
$ws.Cells.Item(5, 10).Value = 0.0246032179869919
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.509196
$ws.Cells.Item(5, 14).Value = 1.527588
$ws.Cells.Item(5, 15).Value = 0.02558190413389134
$ws.Cells.Item(5, 16).Value = 0.02558190413389134
$ws.Cells.Item(5, 17).Value = 0.01818169184
$ws.Cells.Item(5, 18).Value = 0.16363522656
$ws.Cells.Item(5, 19).Value = 0.0006293971639284579
$ws.Cells.Item(5, 20).Value = 0.0006293971639284579

# Row 6: sCs -> Artn/Gfra1 -> FAPs
$ws.Cells.Item(6, 1).Value  = "sCs"
$ws.Cells.Item(6, 2).Value  = "Artn"
$ws.Cells.Item(6, 3).Value  = "Gfra1"
$ws.Cells.Item(6, 4).Value  = "FAPs"
$ws.Cells.Item(6, 5).Value  = 1
$ws.Cells.Item(6, 6).Value  = 0.3333333333333333
$ws.Cells.Item(6, 7).Value  = 0.03570666666666667
$ws.Cells.Item(6, 8).Value  = 0.10712
$ws.Cells.Item(6, 9).Value  = 0.0246032179869919
$ws.Cells.Item(6, 10).Value = 0.0246032179869919
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 17.676258
$ws.Cells.Item(6, 14).Value = 53.028774
$ws.Cells.Item(6, 15).Value = 0.8880516296316739
$ws.Cells.Item(6, 16).Value = 0.8880516296316739
$ws.Cells.Item(6, 17).Value = 0.6311602523200001
$ws.Cells.Item(6, 18).Value = 5.68044227088
$ws.Cells.Item(6, 19).Value = 0.02184892782753147
$ws.Cells.Item(6, 20).Value = 0.02184892782753147

# Row 7: sCs -> Artn/Gfra1 -> sCs
$ws.Cells.Item(7, 1).Value  = "sCs"
$ws.Cells.Item(7, 2).Value  = "Artn"
$ws.Cells.Item(7, 3).Value  = "Gfra1"
$ws.Cells.Item(7, 4).Value  = "sCs"
$ws.Cells.Item(7, 5).Value  = 1
$ws.Cells.Item(7, 6).Value  = 0.3333333333333333
$ws.Cells.Item(7, 7).Value  = 0.03570666666666667
$ws.Cells.Item(7, 8).Value  = 0.10712
$ws.Cells.Item(7, 9).Value  = 0.0246032179869919
$ws.Cells.Item(7, 10).Value = 0.0246032179869919
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.719084666666667
$ws.Cells.Item(7, 14).Value = 5.157254
$ws.Cells.Item(7, 15).Value = 0.0863664662344347
$ws.Cells.Item(7, 16).Value = 0.0863664662344347
$ws.Cells.Item(7, 17).Value = 0.06138278316444445
$ws.Cells.Item(7, 18).Value = 0.5524450484800001
$ws.Cells.Item(7, 19).Value = 0.002124892995531973
$ws.Cells.Item(7, 20).Value = 0.002124892995531973
